$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits (shared-string backed cell values) -----------------
# D12: "Azure DevOps" -> "Configuração da IDE de deploy automatizado / GITHUB"
$ws.Range("D12").Value = "Configuração da IDE de deploy automatizado / GITHUB"

# D15: "Wi-fi, 3G, 4G, 5G" -> "Conexão com a internet"
$ws.Range("D15").Value = "Conexão com a internet"

# D24: "MySQL e PowerPoint" -> "Workbench e PowerPoint"
$ws.Range("D24").Value = "Workbench e PowerPoint"

# --- Row height change (row 5 grew from 29.25 to 46.8) -----------------
$ws.Rows(5).RowHeight = 46.8

# --- Sheet view / window state changes ---------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 100
$win.ScrollRow = 19
$win.ScrollColumn = 2

# New selection: D24 active cell
$ws.Range("D24").Select() | Out-Null
